{"js": "// The document's only paragraph holds two runs of stray placeholder text\n// (\"sd\" and \"fsd\") bracketing a \"_GoBack\" bookmark:\n//   <w:r>sd</w:r><w:bookmarkStart .../><w:bookmarkEnd .../><w:r>fsd</w:r>\n// The fix removes both text runs (the leftover placeholder text that should\n// have been the properly-formatted/mathml space-group content) while\n// leaving the paragraph, its style, and the bookmark untouched.\n\nconst body = context.document.body;\n\n// \"sd\" is ambiguous on its own (it is also a substring of \"fsd\"), but it is\n// the very first match in document order, so resolve it first while the\n// text is still \"sdfsd\" \u2014 this removes exactly the first run.\nlet hits = body.search(\"sd\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (hits.items.length > 0) {\n  hits.items[0].delete();\n  await context.sync();\n}\n\n// The remaining text is now just \"fsd\" (the second run) and matches\n// unambiguously.\nhits = body.search(\"fsd\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nif (hits.items.length > 0) {\n  hits.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# The document's only paragraph holds two runs of stray placeholder text\n# (\"sd\" and \"fsd\") bracketing a \"_GoBack\" bookmark:\n#   <w:r>sd</w:r><w:bookmarkStart .../><w:bookmarkEnd .../><w:r>fsd</w:r>\n# The fix removes both text runs (the leftover placeholder text that should\n# have been the properly-formatted/mathml space-group content) while\n# leaving the paragraph, its style, and the bookmark untouched.\n\n$d = $word.ActiveDocument\n\n# Character offsets are stable across the two deletes only if we remove the\n# later run first, so the first run's [0,2) offsets stay valid.\n$secondRun = $d.Range(2, 5)\nWrite-Output $secondRun.Text\n$secondRun.Delete()\n\n$firstRun = $d.Range(0, 2)\nWrite-Output $firstRun.Text\n$firstRun.Delete()\n"}
